$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) "longform" sheet: append a new data row (row 32) for the
#    2020-03-24 18:00 snapshot, mirroring the layout/style of row 31.
# ---------------------------------------------------------------------
$long = $wb.Worksheets.Item("longform")

# Copy formatting (incl. the date number format) from the row above so the
# new row reuses the existing style instead of creating a new one.
$long.Range("A31:AP31").Copy()
$long.Range("A32:AP32").PasteSpecial(-4122)

$long.Cells.Item(32, 1).Value = "live"
$long.Cells.Item(32, 2).Value = "govt_canada_ph"
$long.Cells.Item(32, 3).Value = 43914.75

$longValues = @{
    4  = 617;  5  = 0;   6  = 13;
    7  = 358;  8  = 0;   9  = 2;
    10 = 72;   11 = 0;   12 = 0;
    13 = 11;   14 = 10;  15 = 0;
    16 = 588;  17 = 0;   18 = 8;
    19 = 221;  20 = 792; 21 = 4;
    22 = 4;    23 = 31;  24 = 0;
    25 = 18;   26 = 0;   27 = 0;
    28 = 51;   29 = 0;   30 = 0;
    31 = 3;    32 = 0;   33 = 0;
    34 = 13;   35 = 0;   36 = 0;
    37 = 2;    38 = 0;   39 = 0;
    40 = 1;    41 = 0;   42 = 0;
}
foreach ($col in $longValues.Keys) {
    $long.Cells.Item(32, $col).Value = $longValues[$col]
}

# Recreate the view state recorded for this sheet after the edit
# (freeze pane stays active; only the selection / scroll position move).
$long.Range("C32").Select()

# ---------------------------------------------------------------------
# 2) "shortform" sheet: append the equivalent long-format rows (327-365)
#    for the same 2020-03-24 18:00 snapshot: one row per
#    province x case_type combination.
# ---------------------------------------------------------------------
$short = $wb.Worksheets.Item("shortform")

$short.Range("A326:F326").Copy()
$short.Range("A327:F365").PasteSpecial(-4122)

$shortRows = @(
    @("BC",    "conf",   617),
    @("BC",    "prob",   0),
    @("BC",    "deaths", 13),
    @("AB",    "conf",   358),
    @("AB",    "prob",   0),
    @("AB",    "deaths", 2),
    @("SK",    "conf",   72),
    @("SK",    "prob",   0),
    @("SK",    "deaths", 0),
    @("MB",    "conf",   11),
    @("MB",    "prob",   10),
    @("MB",    "deaths", 0),
    @("ON",    "conf",   588),
    @("ON",    "prob",   0),
    @("ON",    "deaths", 8),
    @("QC",    "conf",   221),
    @("QC",    "prob",   792),
    @("QC",    "deaths", 4),
    @("NL",    "conf",   4),
    @("NL",    "prob",   31),
    @("NL",    "deaths", 0),
    @("NB",    "conf",   18),
    @("NB",    "prob",   0),
    @("NB",    "deaths", 0),
    @("NS",    "conf",   51),
    @("NS",    "prob",   0),
    @("NS",    "deaths", 0),
    @("PEI",   "conf",   3),
    @("PEI",   "prob",   0),
    @("PEI",   "deaths", 0),
    @("Repat", "conf",   13),
    @("Repat", "prob",   0),
    @("Repat", "deaths", 0),
    @("YK",    "conf",   2),
    @("YK",    "prob",   0),
    @("YK",    "deaths", 0),
    @("NT",    "conf",   1),
    @("NT",    "prob",   0),
    @("NT",    "deaths", 0)
)

$row = 327
foreach ($entry in $shortRows) {
    $short.Cells.Item($row, 1).Value = "live"
    $short.Cells.Item($row, 2).Value = "govt_canada_ph"
    $short.Cells.Item($row, 3).Value = 43914.75
    $short.Cells.Item($row, 4).Value = $entry[0]
    $short.Cells.Item($row, 5).Value = $entry[1]
    $short.Cells.Item($row, 6).Value = $entry[2]
    $row = $row + 1
}

# Recreate the view state recorded for this sheet after the edit.
$short.Range("A327:C365").Select()
